$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.41903
$ws.Range("H2").Value = 1.25709
$ws.Range("M2").Value = 62.87391666666667
$ws.Range("N2").Value = 188.62175
$ws.Range("O2").Value = 0.5484251561826182
$ws.Range("P2").Value = 0.5484251561826182
$ws.Range("Q2").Value = 26.34605730083333
$ws.Range("R2").Value = 237.1145157074999
$ws.Range("S2").Value = 0.5484251561826182
$ws.Range("T2").Value = 0.5484251561826182

# Row 3
$ws.Range("G3").Value = 0.41903
$ws.Range("H3").Value = 1.25709
$ws.Range("O3").Value = 0.3074331251635
$ws.Range("P3").Value = 0.3074331251635
$ws.Range("Q3").Value = 14.76892633465333
$ws.Range("R3").Value = 132.92033701188
$ws.Range("S3").Value = 0.3074331251635
$ws.Range("T3").Value = 0.3074331251635

# Row 4
$ws.Range("G4").Value = 0.41903
$ws.Range("H4").Value = 1.25709
$ws.Range("M4").Value = 7.926563000000001
$ws.Range("N4").Value = 23.779689
$ws.Range("O4").Value = 0.06914038096772557
$ws.Range("P4").Value = 0.06914038096772555
$ws.Range("Q4").Value = 3.32146769389
$ws.Range("R4").Value = 29.89320924501
$ws.Range("S4").Value = 0.06914038096772557
$ws.Range("T4").Value = 0.06914038096772555

# Row 5
$ws.Range("G5").Value = 0.41903
$ws.Range("H5").Value = 1.25709
$ws.Range("M5").Value = 8.598489333333333
$ws.Range("N5").Value = 25.795468
$ws.Range("O5").Value = 0.07500133768615619
$ws.Range("P5").Value = 0.07500133768615617
$ws.Range("Q5").Value = 3.603024985346666
$ws.Range("R5").Value = 32.42722486811999
$ws.Range("S5").Value = 0.07500133768615619
$ws.Range("T5").Value = 0.07500133768615617
